$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header fixes ---
# B1 "locaclizacion" -> "locacalizacion" (typo fix, same column/position)
$ws.Range("B1").Value = "locacalizacion"

# --- Row 2 data fixes: Juan -> jorge ---
$ws.Range("A2").Value = "jorge"

# B2 was a time value (0.53138888888888891) with style s=5 (numFmtId 21, h:mm AM/PM).
# It becomes the literal text "18:13:14:12S" and the style itself changes:
# new font (Courier New, 10pt, color FF6A8759), General number format, vertically centered.
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("B2").Font.Name = "Courier New"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").Font.Color = 5867370
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").NumberFormat = "general"

# C2 keeps its hyperlink (rels untouched) but display text changes
$ws.Range("C2").Value = "jorge@email.es"

# D2 id value change
$ws.Range("D2").Value = "ID4"

# --- Selection change ---
$ws.Range("A1:E2").Select()
